$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 43 (shifts rows 43:63 down to 44:64)
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with the teacher's name field
$ws.Range("A43").Value = "姓名"
$ws.Range("B43").Value = "tname"
$ws.Range("C43").Value = "varchar(50)"
$ws.Range("D43").Value = 1

# Update the selection to match the recorded state after the edit
$ws.Range("G44").Select()
